$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  111.193776,
  117.66211199999999,
  113.349887999999,
  115.506,
  112.88786399999999,
  114.27393600000001,
  117.970128,
  113.19588,
  113.19588,
  113.349887999999,
  113.65790399999899,
  109.49968800000001,
  114.581952,
  113.19588,
  116.584056,
  113.96592,
  113.503896,
  117.354096,
  113.349887999999,
  113.19588,
  113.503896,
  113.65790399999899,
  113.811911999999,
  114.427943999999,
  117.66211199999999,
  113.503896,
  113.349887999999,
  112.88786399999999,
  109.807704,
  113.811911999999
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("E10").Select()
